# Add data for 2022-12-19
# Updates the "through" date in the sheet name and header label from
# December 10 to December 11, and bumps/adds carjacking counts for the
# newly-included day across several neighborhood rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Through 2022-12-11"

# Update the header cell text (B1, shared string) describing the partial month
$ws.Cells.Item(1, 2).Value = "December 2022 (through December 11)"

# Cell value updates: [row, column, newValue]
$updates = @(
    @(2, 62, 4),    # BJ2: 3 -> 4
    @(7, 2, 3),      # B7: 1 -> 3
    @(7, 26, 3),     # Z7: 1 -> 3
    @(19, 14, 2),    # N19: 1 -> 2
    @(20, 26, 7),    # Z20: 6 -> 7
    @(21, 2, 1),     # B21: new -> 1
    @(22, 62, 1),    # BJ22: new -> 1
    @(28, 14, 2),    # N28: 1 -> 2
    @(32, 14, 1),    # N32: new -> 1
    @(32, 26, 2),    # Z32: 1 -> 2
    @(33, 38, 2),    # AL33: 1 -> 2
    @(35, 14, 2),    # N35: 1 -> 2
    @(36, 2, 1),     # B36: new -> 1
    @(36, 14, 2),    # N36: 1 -> 2
    @(36, 38, 1),    # AL36: new -> 1
    @(36, 86, 2),    # CH36: 1 -> 2
    @(37, 2, 1),     # B37: new -> 1
    @(39, 74, 1),    # BV39: new -> 1
    @(40, 14, 2),    # N40: 1 -> 2
    @(42, 26, 1),    # Z42: new -> 1
    @(45, 2, 2),     # B45: 1 -> 2
    @(49, 14, 1),    # N49: new -> 1
    @(53, 14, 1),    # N53: new -> 1
    @(57, 74, 1),    # BV57: new -> 1
    @(64, 26, 2),    # Z64: 1 -> 2
    @(72, 26, 1),    # Z72: new -> 1
    @(84, 14, 1),    # N84: new -> 1
    @(84, 62, 1)     # BJ84: new -> 1
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $ws.Cells.Item($row, $col).Value = $val
}
